$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10: September (through 09-04) -> (through 09-05), with new counts
$ws.Range("A10").Value = "September (through 09-05)"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 18
$ws.Range("H10").Value = 22

# Update row 11 (Total) with new cumulative counts
$ws.Range("B11").Value = 197
$ws.Range("C11").Value = 390
$ws.Range("D11").Value = 565
$ws.Range("E11").Value = 497
$ws.Range("F11").Value = 360
$ws.Range("G11").Value = 802
$ws.Range("H11").Value = 1093
